# theodoiphongmay.xlsx - "them dung vi tri cua cac truong"
$wb = $excel.ActiveWorkbook

# 1) Rename sheet "Tuần 2" -> "Tuần 3"
$wsWeek = $wb.Worksheets.Item("Tuần 2")
$wsWeek.Name = "Tuần 3"

# 2) Fix header labels ("Thứ 2.3" -> "Thứ 2", "Thứ 3.2" -> "Thứ 3")
$wsWeek.Range("A3").Value = "Thứ 2"
$wsWeek.Range("A13").Value = "Thứ 3"

# 3) Fill in the "Thứ 5" block (rows 33-41) with the new class/lesson/teacher info.
#    Column C holds a class code like "1/1" that looks numeric, so it is entered with a
#    leading apostrophe - exactly how Excel records a manually quote-prefixed text entry
#    (produces the quotePrefix="1" cell style).
$wsWeek.Range("C33").Value = "'1/1"
$wsWeek.Range("E33").Value = "test vi tri v3"
$wsWeek.Range("G33").Value = "Cô Nguyên"

$wsWeek.Range("C35").Value = "'1/2"
$wsWeek.Range("E35").Value = "test vi tri v3"
$wsWeek.Range("G35").Value = "Cô Nguyên"

$wsWeek.Range("C37").Value = "'1/3"
$wsWeek.Range("E37").Value = "test vi tri v3"
$wsWeek.Range("G37").Value = "Cô Nguyên"

$wsWeek.Range("C39").Value = "'1/2"
$wsWeek.Range("E39").Value = "test vi tri v3"
$wsWeek.Range("G39").Value = "Cô Nguyên"

$wsWeek.Range("C41").Value = "'1/3"
$wsWeek.Range("E41").Value = "test vi tri v3"
$wsWeek.Range("G41").Value = "Cô Nguyên"

# 4) Row height: every row now uses the (new) default height of 15.75 instead of an
#    explicit per-row override.
$wsWeek.Rows("1:55").RowHeight = 15.75

$wsOther = $wb.Worksheets.Item("Trang_tính1")
$wsOther.Rows("1:55").RowHeight = 15.75
